$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.113.55"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.319.28"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'313.51"
$ws.Range("E5").Value = "  -5.67%  "
$ws.Range("D6").Value = "'106.17"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -2.39%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").Value = "'40.48"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").Value = "'0.0916"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "'8.30"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "'0.984"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'15.65"
$ws.Range("E15").Value = "  -5.24%  "
$ws.Range("D16").Value = "2.668.17"
$ws.Range("D17").Value = "2.315.69"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "42.145.69"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "'7.71"
$ws.Range("E19").Value = "  -4.70%  "
$ws.Range("D21").Value = "'74.82"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "'3.49"
$ws.Range("E22").Value = "  -7.83%  "
$ws.Range("D23").Value = "'260.93"
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("D24").Value = "'2.31"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'9.32"
$ws.Range("E25").Value = "  -8.95%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'22.85"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'35.77"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "'0.0900"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'163.08"
$ws.Range("E32").Value = "  -7.57%  "
$ws.Range("E33").Value = "  -5.83%  "
$ws.Range("D34").Value = "'5.88"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.118"
$ws.Range("E36").Value = "  +11.82%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").Value = "'0.0354"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "  -6.31%  "
$ws.Range("E40").Value = "  -4.51%  "
$ws.Range("D41").Value = "'98.32"
$ws.Range("E41").Value = "  +6.27%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'71.24"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.47"
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("D44").Value = "'0.231"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'12.31"
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("D47").Value = "'112.10"
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "'9.02"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").Value = "'74.77"
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("E51").Value = "  -0.71%  "
